# Add a new worksheet "k_p_z_min" after the "k_p" sheet, populate it with
# the P / k_p_z_min table, and make it the active/selected sheet.

$wb = $excel.ActiveWorkbook

# Add the new sheet right after the "k_p" sheet.
$kpSheet = $wb.Worksheets.Item("k_p")
$newSheet = $wb.Worksheets.Add($null, $kpSheet)
$newSheet.Name = "k_p_z_min"

# Header row.
$newSheet.Range("A1").Value = "P"
$newSheet.Range("B1").Value = "k_p_z_min"

# Data rows.
$data = @(
    @(500, 0.08),
    @(1000, 0.1),
    @(1500, 0.12),
    @(2000, 0.14),
    @(2500, 0.15)
)

$row = 2
foreach ($pair in $data) {
    $newSheet.Cells.Item($row, 1).Value = $pair[0]
    $newSheet.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}

# Select B2 on the new sheet and make it the active sheet/tab.
$newSheet.Range("B2").Select()
$newSheet.Activate()
